$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns Q = "PD", R = "N2"
$ws.Cells.Item(1, 17).Value = "PD"
$ws.Cells.Item(1, 18).Value = "N2"

# Data rows 2-79: column Q (PD polygon code) and column R (N2 zone description)
$ws.Cells.Item(2, 17).Value = "RET-B"
$ws.Cells.Item(2, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(3, 17).Value = "AGU-B"
$ws.Cells.Item(3, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(4, 17).Value = "DEV-D"
$ws.Cells.Item(4, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(5, 17).Value = "COG-O"
$ws.Cells.Item(5, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(6, 17).Value = "BLO-R"
$ws.Cells.Item(6, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(7, 17).Value = "COG-G"
$ws.Cells.Item(7, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(8, 17).Value = "ATH-N"
$ws.Cells.Item(8, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(9, 17).Value = "VCR-E"
$ws.Cells.Item(9, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(10, 17).Value = "AGU-B"
$ws.Cells.Item(10, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(11, 17).Value = "CEN-G"
$ws.Cells.Item(11, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(12, 17).Value = "NRA-M"
$ws.Cells.Item(12, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(13, 17).Value = "ALM-C"
$ws.Cells.Item(13, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(14, 17).Value = "CEN-C"
$ws.Cells.Item(14, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(15, 17).Value = "DEV-L"
$ws.Cells.Item(15, 18).Value = "ARATO-25058.PO.2DEV"
$ws.Cells.Item(16, 17).Value = "COG-G"
$ws.Cells.Item(16, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(17, 17).Value = "VCR-H"
$ws.Cells.Item(17, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(18, 17).Value = "ATH-?"
$ws.Cells.Item(18, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(19, 17).Value = "CEN-G"
$ws.Cells.Item(19, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(20, 17).Value = "NRA-R"
$ws.Cells.Item(20, 18).Value = "ARATO-25058.PO.2NRA"
$ws.Cells.Item(21, 17).Value = "ALM-N"
$ws.Cells.Item(21, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(22, 17).Value = "PAV-N"
$ws.Cells.Item(22, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(23, 17).Value = "BLO-G"
$ws.Cells.Item(23, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(24, 17).Value = "ATH-D"
$ws.Cells.Item(24, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(25, 17).Value = "VCR-D"
$ws.Cells.Item(25, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(26, 17).Value = "RET-D"
$ws.Cells.Item(26, 18).Value = "ARATO-25058.PO.1RET"
$ws.Cells.Item(27, 17).Value = "ALM-F"
$ws.Cells.Item(27, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(28, 17).Value = "CON-C"
$ws.Cells.Item(28, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(29, 17).Value = "CEN-M"
$ws.Cells.Item(29, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(30, 17).Value = "CEN-B"
$ws.Cells.Item(30, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(31, 17).Value = "CEN-E"
$ws.Cells.Item(31, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(32, 17).Value = "RET-H"
$ws.Cells.Item(32, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(33, 17).Value = "COG-D"
$ws.Cells.Item(33, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(34, 17).Value = "COG-K"
$ws.Cells.Item(34, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(35, 17).Value = "VCR-M"
$ws.Cells.Item(35, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(36, 17).Value = "ATH-P"
$ws.Cells.Item(36, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(37, 17).Value = "ATH-P"
$ws.Cells.Item(37, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(38, 17).Value = "CEN-G"
$ws.Cells.Item(38, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(39, 17).Value = "CLI-M"
$ws.Cells.Item(39, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(40, 17).Value = "BLO-S"
$ws.Cells.Item(40, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(41, 17).Value = "CEN-N"
$ws.Cells.Item(41, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(42, 17).Value = "CLI-D"
$ws.Cells.Item(42, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(43, 17).Value = "CLI-F"
$ws.Cells.Item(43, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(44, 17).Value = "RET-R"
$ws.Cells.Item(44, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(45, 17).Value = "VCR-J"
$ws.Cells.Item(45, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(46, 17).Value = "BLO-G"
$ws.Cells.Item(46, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(47, 17).Value = "ATH-P"
$ws.Cells.Item(47, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(48, 17).Value = "CLI-O"
$ws.Cells.Item(48, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(49, 17).Value = "ATH-I"
$ws.Cells.Item(49, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(50, 17).Value = "CEN-E"
$ws.Cells.Item(50, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(51, 17).Value = "DEV-L"
$ws.Cells.Item(51, 18).Value = "ARATO-25058.PO.2DEV"
$ws.Cells.Item(52, 17).Value = "CLI-D"
$ws.Cells.Item(52, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(53, 17).Value = "RET-E"
$ws.Cells.Item(53, 18).Value = "ARATO-25058.PO.1RET"
$ws.Cells.Item(54, 17).Value = "AGU-F"
$ws.Cells.Item(54, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(55, 17).Value = "AGU-L"
$ws.Cells.Item(55, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(56, 17).Value = "VCR-O"
$ws.Cells.Item(56, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(57, 17).Value = "ALM-H"
$ws.Cells.Item(57, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(58, 17).Value = "CLI-B"
$ws.Cells.Item(58, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(59, 17).Value = "CEN-C"
$ws.Cells.Item(59, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(60, 17).Value = "CLI-M"
$ws.Cells.Item(60, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(61, 17).Value = "PCH-M"
$ws.Cells.Item(61, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(62, 17).Value = "ALM-O"
$ws.Cells.Item(62, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(63, 17).Value = "PCH-F"
$ws.Cells.Item(63, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(64, 17).Value = "PCH-G"
$ws.Cells.Item(64, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(65, 17).Value = "NRA-I"
$ws.Cells.Item(65, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(66, 17).Value = "COG-C"
$ws.Cells.Item(66, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(67, 17).Value = "ATH-J"
$ws.Cells.Item(67, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(68, 17).Value = "ALM-O"
$ws.Cells.Item(68, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(69, 17).Value = "VCR-D"
$ws.Cells.Item(69, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(70, 17).Value = "ATH-I"
$ws.Cells.Item(70, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(71, 17).Value = "CON-K"
$ws.Cells.Item(71, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(72, 17).Value = "ALM-O"
$ws.Cells.Item(72, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(73, 17).Value = "COG-M"
$ws.Cells.Item(73, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(74, 17).Value = "COG-M"
$ws.Cells.Item(74, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(75, 17).Value = "COG-L"
$ws.Cells.Item(75, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(76, 17).Value = "ATH-B"
$ws.Cells.Item(76, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(77, 17).Value = "CLI-F"
$ws.Cells.Item(77, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(78, 17).Value = "CON-H"
$ws.Cells.Item(78, 18).Value = "Fuera de Poligono OVL"
$ws.Cells.Item(79, 17).Value = "CON-H"
$ws.Cells.Item(79, 18).Value = "Fuera de Poligono OVL"

# Match header styling (bold, centered, bordered) used by the existing header cells
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
